# Fixing model-traits.cfg.* files to add entries for model_label, label (trait),
# and label_short (abbreviated descriptor for trait).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("trait" and everything to its right
# shifts one column further right, e.g. trait: B->C, filter: C->D, etc.).
$ws.Columns("B:B").Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "model_label"

# Rows 2-22 are the per-year models (2011/2012/2013); model_label simply
# mirrors the year already stored in column A. Set the shared range B3:B22
# first, then the lone leading cell B2, so the shared-formula grouping lines
# up the same way the existing columns on the sheet are grouped.
$ws.Range("B3:B22").Formula = "=A3"
$ws.Range("B2").Formula = "=A2"

# Rows 23-29 correspond to the combined "all-years" model; label them
# accordingly.
$ws.Range("B23:B29").Value = "All Years"

# The "trait~1" / "trait~year" helper formulas in column E (previously D,
# before the new column was inserted) now need to read from column C
# (previously B, where "trait" now lives). Re-assert them as range formulas
# so Excel keeps/re-builds the shared-formula grouping cleanly.
$ws.Range("E3:E22").Formula = "=_xlfn.CONCAT(C3,""~1"")"
$ws.Range("E24:E29").Formula = "=_xlfn.CONCAT(C24,""~year"")"

# Move the active selection, matching the edited workbook.
$ws.Range("B32").Select()
